{"js": "// Update the worksheet date and the twenty-five multiplication problems.\nconst replacements = [\n    { old: \"2026-02-20 Friday\", new: \"2026-02-21 Saturday\" },\n    { old: \"375\u00d78=\", new: \"847\u00d79=\" },\n    { old: \"252\u00d78=\", new: \"116\u00d73=\" },\n    { old: \"644\u00d76=\", new: \"320\u00d74=\" },\n    { old: \"845\u00d78=\", new: \"706\u00d76=\" },\n    { old: \"303\u00d79=\", new: \"994\u00d76=\" },\n    { old: \"813\u00d77=\", new: \"184\u00d73=\" },\n    { old: \"233\u00d76=\", new: \"990\u00d75=\" },\n    { old: \"853\u00d76=\", new: \"833\u00d79=\" },\n    { old: \"764\u00d75=\", new: \"838\u00d74=\" },\n    { old: \"872\u00d74=\", new: \"289\u00d73=\" },\n    { old: \"854\u00d77=\", new: \"535\u00d74=\" },\n    { old: \"920\u00d73=\", new: \"112\u00d78=\" },\n    { old: \"420\u00d72=\", new: \"529\u00d79=\" },\n    { old: \"620\u00d74=\", new: \"998\u00d72=\" },\n    { old: \"915\u00d78=\", new: \"723\u00d75=\" },\n    { old: \"530\u00d78=\", new: \"609\u00d76=\" },\n    { old: \"686\u00d76=\", new: \"671\u00d76=\" },\n    { old: \"877\u00d74=\", new: \"401\u00d79=\" },\n    { old: \"341\u00d77=\", new: \"647\u00d74=\" },\n    { old: \"498\u00d77=\", new: \"914\u00d77=\" },\n    { old: \"169\u00d73=\", new: \"701\u00d79=\" },\n    { old: \"618\u00d72=\", new: \"612\u00d72=\" },\n    { old: \"841\u00d74=\", new: \"732\u00d73=\" },\n    { old: \"472\u00d77=\", new: \"417\u00d73=\" },\n    { old: \"293\u00d76=\", new: \"878\u00d79=\" }\n];\n\nconst body = context.document.body;\n\nfor (const r of replacements) {\n    const results = body.search(r.old, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < results.items.length; i++) {\n        results.items[i].insertText(r.new, Word.InsertLocation.replace);\n    }\n    await context.sync();\n}\n", "ps1": "# Update the worksheet date and the twenty-five multiplication problems.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2026-02-20 Friday\"; New = \"2026-02-21 Saturday\" },\n    @{ Old = \"375\u00d78=\";            New = \"847\u00d79=\" },\n    @{ Old = \"252\u00d78=\";            New = \"116\u00d73=\" },\n    @{ Old = \"644\u00d76=\";            New = \"320\u00d74=\" },\n    @{ Old = \"845\u00d78=\";            New = \"706\u00d76=\" },\n    @{ Old = \"303\u00d79=\";            New = \"994\u00d76=\" },\n    @{ Old = \"813\u00d77=\";            New = \"184\u00d73=\" },\n    @{ Old = \"233\u00d76=\";            New = \"990\u00d75=\" },\n    @{ Old = \"853\u00d76=\";            New = \"833\u00d79=\" },\n    @{ Old = \"764\u00d75=\";            New = \"838\u00d74=\" },\n    @{ Old = \"872\u00d74=\";            New = \"289\u00d73=\" },\n    @{ Old = \"854\u00d77=\";            New = \"535\u00d74=\" },\n    @{ Old = \"920\u00d73=\";            New = \"112\u00d78=\" },\n    @{ Old = \"420\u00d72=\";            New = \"529\u00d79=\" },\n    @{ Old = \"620\u00d74=\";            New = \"998\u00d72=\" },\n    @{ Old = \"915\u00d78=\";            New = \"723\u00d75=\" },\n    @{ Old = \"530\u00d78=\";            New = \"609\u00d76=\" },\n    @{ Old = \"686\u00d76=\";            New = \"671\u00d76=\" },\n    @{ Old = \"877\u00d74=\";            New = \"401\u00d79=\" },\n    @{ Old = \"341\u00d77=\";            New = \"647\u00d74=\" },\n    @{ Old = \"498\u00d77=\";            New = \"914\u00d77=\" },\n    @{ Old = \"169\u00d73=\";            New = \"701\u00d79=\" },\n    @{ Old = \"618\u00d72=\";            New = \"612\u00d72=\" },\n    @{ Old = \"841\u00d74=\";            New = \"732\u00d73=\" },\n    @{ Old = \"472\u00d77=\";            New = \"417\u00d73=\" },\n    @{ Old = \"293\u00d76=\";            New = \"878\u00d79=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute([ref]$r.Old, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$r.New, 2)\n}\n"}
